# Apply the diff: update the "as of" date in the confidentiality banner
# (2021-04-29 -> 2021-04-30) and refresh the Weight / Percent Change
# figures for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected; unlock it so the cells can be edited, then
# restore protection afterwards.
$ws.Unprotect()

# --- Disclaimer banner text (shared string used by A7) ---------------
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-30 for illustrative purposes only and are subject to change."
# Re-fit the row height so the newline in the text above doesn't leave a
# stray custom row height behind.
$ws.Rows("7").AutoFit()

# --- Weight (column D) and Percent Change (column E) figures ---------
$ws.Range("D2").Value = 0.8437033349925336
$ws.Range("E2").Value = -0.01251580278128939

$ws.Range("D3").Value = 0.1562966650074664
$ws.Range("E3").Value = -0.01765241128298467

$ws.Range("E4").Value = -0.01331863755955331

# Restore sheet protection.
$ws.Protect()
